# Update the account credentials on the "Registration" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

$ws.Range("B2").Value = "test44332@yopmail.com"
$ws.Range("B3").Value = "C!0ud24@h2Ah"
